# Scheduled data-refresh: update market-price-derived columns (H, I, J, K,
# L, M, N -> currentAveragePrice*/LevePrice*/LeveProfit*) across the
# ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets with freshly pulled prices.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 419806.12
$ws.Range("J17").Value = 419806.12
$ws.Range("L17").Value = 1259418.36
$ws.Range("N17").Value = -1259754.36

$ws.Range("I43").Value = 2999
$ws.Range("J43").Value = 4000
$ws.Range("K43").Value = 2999
$ws.Range("L43").Value = 4000
$ws.Range("M43").Value = -2930
$ws.Range("N43").Value = -4138

$ws.Range("H58").Value = 1000.8571
$ws.Range("I58").Value = 100
$ws.Range("J58").Value = 1361.2
$ws.Range("K58").Value = 300
$ws.Range("L58").Value = 4083.6
$ws.Range("M58").Value = -150
$ws.Range("N58").Value = -4383.6

$ws.Range("H112").Value = 5810189
$ws.Range("J112").Value = 7746251
$ws.Range("L112").Value = 23238753
$ws.Range("N112").Value = -23240969

$ws.Range("H138").Value = 3505.1
$ws.Range("I138").Value = 1987.28
$ws.Range("J138").Value = 4011.04
$ws.Range("K138").Value = 5961.84
$ws.Range("L138").Value = 12033.12
$ws.Range("M138").Value = -821.8400000000001
$ws.Range("N138").Value = -22313.12

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 11617.571
$ws.Range("I28").Value = 11617.571
$ws.Range("K28").Value = 11617.571
$ws.Range("M28").Value = -11425.571

$ws.Range("H74").Value = 70677.62
$ws.Range("I74").Value = 70677.62
$ws.Range("K74").Value = 70677.62
$ws.Range("M74").Value = -69803.62

$ws.Range("H77").Value = 70677.62
$ws.Range("I77").Value = 70677.62
$ws.Range("K77").Value = 353388.1
$ws.Range("M77").Value = -349020.1

$ws.Range("H99").Value = 11617.571
$ws.Range("I99").Value = 11617.571
$ws.Range("K99").Value = 11617.571
$ws.Range("M99").Value = -8622.571

$ws.Range("H110").Value = 7311.3
$ws.Range("I110").Value = 6041.6
$ws.Range("K110").Value = 6041.6
$ws.Range("M110").Value = -3996.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2035.8
$ws.Range("I20").Value = 1412.0333
$ws.Range("K20").Value = 1412.0333
$ws.Range("M20").Value = -1165.0333

$ws.Range("H99").Value = 7942.0415
$ws.Range("I99").Value = 9890.6
$ws.Range("J99").Value = 4694.4443
$ws.Range("K99").Value = 9890.6
$ws.Range("L99").Value = 4694.4443
$ws.Range("M99").Value = -8392.6
$ws.Range("N99").Value = -7690.4443

$ws.Range("H105").Value = 2058.1667
$ws.Range("I105").Value = 2127.182
$ws.Range("J105").Value = 1299
$ws.Range("K105").Value = 2127.182
$ws.Range("L105").Value = 1299
$ws.Range("M105").Value = -380.1819999999998
$ws.Range("N105").Value = -4793

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H37").Value = 29000
$ws.Range("J37").Value = 29000
$ws.Range("L37").Value = 29000
$ws.Range("N37").Value = -29214

$ws.Range("H68").Value = 420000
$ws.Range("J68").Value = 420000
$ws.Range("L68").Value = 420000
$ws.Range("N68").Value = -421498

$ws.Range("H71").Value = 420000
$ws.Range("J71").Value = 420000
$ws.Range("L71").Value = 1260000
$ws.Range("N71").Value = -1267488

$ws.Range("H75").Value = 30000
$ws.Range("J75").Value = 30000
$ws.Range("L75").Value = 30000
$ws.Range("N75").Value = -31996

$ws.Range("H78").Value = 30000
$ws.Range("J78").Value = 30000
$ws.Range("L78").Value = 90000
$ws.Range("N78").Value = -99984

$ws.Range("H81").Value = 44249.25
$ws.Range("J81").Value = 43999.332
$ws.Range("L81").Value = 43999.332
$ws.Range("N81").Value = -45995.332

$ws.Range("H84").Value = 44249.25
$ws.Range("J84").Value = 43999.332
$ws.Range("L84").Value = 131997.996
$ws.Range("N84").Value = -141981.996

$ws.Range("H120").Value = 48304.555
$ws.Range("I120").Value = 48900
$ws.Range("J120").Value = 48134.43
$ws.Range("K120").Value = 48900
$ws.Range("L120").Value = 48134.43
$ws.Range("M120").Value = -45271
$ws.Range("N120").Value = -55392.43

$ws.Range("H121").Value = 47081
$ws.Range("J121").Value = 46774.668
$ws.Range("L121").Value = 46774.668
$ws.Range("N121").Value = -49394.668

$ws.Range("H122").Value = 1445.6
$ws.Range("I122").Value = 1178.4615
$ws.Range("J122").Value = 3182
$ws.Range("K122").Value = 3535.3845
$ws.Range("L122").Value = 9546
$ws.Range("M122").Value = -1085.3845
$ws.Range("N122").Value = -14446

$ws.Range("H134").Value = 5930.4165
$ws.Range("I134").Value = 6195.636
$ws.Range("J134").Value = 3013
$ws.Range("K134").Value = 18586.908
$ws.Range("L134").Value = 9039
$ws.Range("M134").Value = -16051.908
$ws.Range("N134").Value = -14109

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 520.53845
$ws.Range("I86").Value = 586.1111
$ws.Range("K86").Value = 1758.3333
$ws.Range("M86").Value = -572.3332999999998

$ws.Range("H89").Value = 520.53845
$ws.Range("I89").Value = 586.1111
$ws.Range("K89").Value = 5274.9999
$ws.Range("M89").Value = 653.0001000000002

$ws.Range("H98").Value = 317.4
$ws.Range("I98").Value = 317.4
$ws.Range("K98").Value = 952.1999999999999
$ws.Range("M98").Value = 545.8000000000001

$ws.Range("H113").Value = 1855.4615
$ws.Range("J113").Value = 2100.6
$ws.Range("L113").Value = 6301.799999999999
$ws.Range("N113").Value = -10641.8

$ws.Range("H120").Value = 2900
$ws.Range("I120").Value = 2900
$ws.Range("K120").Value = 8700
$ws.Range("M120").Value = -3862

$ws.Range("H129").Value = 2190.1667
$ws.Range("I129").Value = 1909.8572
$ws.Range("K129").Value = 5729.571599999999
$ws.Range("M129").Value = -729.5715999999993

$ws.Range("H132").Value = 5407.4
$ws.Range("I132").Value = 6227.2085
$ws.Range("J132").Value = 2128.1667
$ws.Range("K132").Value = 56044.8765
$ws.Range("L132").Value = 19153.5003
$ws.Range("M132").Value = -53514.8765
$ws.Range("N132").Value = -24213.5003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 6000
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 6000
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 6000
$ws.Range("M97").ClearContents()
$ws.Range("N97").Value = -6992

$ws.Range("H102").Value = 2351.0312
$ws.Range("I102").Value = 1049.4333
$ws.Range("K102").Value = 1049.4333
$ws.Range("M102").Value = 572.5667000000001

$ws.Range("H122").Value = 1836.6842
$ws.Range("I122").Value = 1301.1538
$ws.Range("J122").Value = 2997
$ws.Range("K122").Value = 3903.4614
$ws.Range("L122").Value = 8991
$ws.Range("M122").Value = -1453.4614
$ws.Range("N122").Value = -13891

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5684.52
$ws.Range("J7").Value = 5774.625
$ws.Range("L7").Value = 5774.625
$ws.Range("N7").Value = -5998.625

$ws.Range("H40").Value = 3767.5925
$ws.Range("I40").Value = 3643.2693
$ws.Range("J40").Value = 7000
$ws.Range("K40").Value = 3643.2693
$ws.Range("L40").Value = 7000
$ws.Range("M40").Value = -3507.2693
$ws.Range("N40").Value = -7272

$ws.Range("H100").Value = 4829
$ws.Range("I100").Value = 4499
$ws.Range("K100").Value = 4499
$ws.Range("M100").Value = -3958

$ws.Range("H126").Value = 5684.52
$ws.Range("J126").Value = 5774.625
$ws.Range("L126").Value = 17323.875
$ws.Range("N126").Value = -22263.875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H93").Value = 59993
$ws.Range("I93").Value = 59993
$ws.Range("K93").Value = 59993
$ws.Range("M93").Value = -57497
